$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J2").Value = "{'spatial': {'bbox': [[0, -89, 360, 89]]}, 'temporal': {'interval': [['2024-03-01T00:00:00Z', '2025-12-06T00:00:00Z']]}}"
$ws.Range("J6").Value = "{'spatial': {'bbox': [[-25, 30, 45, 72]]}, 'temporal': {'interval': [['2024-01-17T00:00:00Z', '2025-12-06T00:00:00Z']]}}"
$ws.Range("J8").Value = "{'spatial': {'bbox': [[-180, -90, 180, 90]]}, 'temporal': {'interval': [['2003-01-01T00:00:00Z', '2025-12-05T00:00:00Z']]}}"
$ws.Range("J11").Value = "{'spatial': {'bbox': [[-180, -90, 180, 90]]}, 'temporal': {'interval': [['2015-01-01T00:00:00Z', '2025-12-06T00:00:00Z']]}}"
$ws.Range("J16").Value = "{'spatial': {'bbox': [[-25, 30, 45, 72]]}, 'temporal': {'interval': [['2022-12-01T00:00:00Z', '2025-12-06T00:00:00Z']]}}"
$ws.Range("J17").Value = "{'spatial': {'bbox': [[0, -89, 360, 89]]}, 'temporal': {'interval': [['2004-01-01T00:00:00Z', '2025-12-05T00:00:00Z']]}}"

$wb.Save()
